$d = $word.ActiveDocument

$replacements = @(
    @{ old = "741÷7=105, 6"; new = "506÷9=56, 2" },
    @{ old = "526÷3=175, 1"; new = "759÷3=253, 0" },
    @{ old = "844÷9=93, 7"; new = "444÷2=222, 0" },
    @{ old = "245÷7=35, 0"; new = "423÷4=105, 3" },
    @{ old = "562÷6=93, 4"; new = "704÷5=140, 4" },
    @{ old = "764÷3=254, 2"; new = "102÷9=11, 3" },
    @{ old = "137÷2=68, 1"; new = "975÷4=243, 3" },
    @{ old = "511÷6=85, 1"; new = "243÷9=27, 0" },
    @{ old = "825÷5=165, 0"; new = "132÷6=22, 0" },
    @{ old = "812÷5=162, 2"; new = "576÷4=144, 0" },
    @{ old = "806÷9=89, 5"; new = "508÷4=127, 0" },
    @{ old = "646÷2=323, 0"; new = "504÷4=126, 0" },
    @{ old = "523÷4=130, 3"; new = "864÷5=172, 4" },
    @{ old = "348÷9=38, 6"; new = "610÷6=101, 4" },
    @{ old = "248÷5=49, 3"; new = "312÷9=34, 6" },
    @{ old = "135÷3=45, 0"; new = "916÷9=101, 7" },
    @{ old = "726÷9=80, 6"; new = "665÷8=83, 1" },
    @{ old = "948÷3=316, 0"; new = "652÷9=72, 4" },
    @{ old = "767÷6=127, 5"; new = "185÷2=92, 1" },
    @{ old = "723÷5=144, 3"; new = "395÷6=65, 5" },
    @{ old = "823÷8=102, 7"; new = "326÷2=163, 0" },
    @{ old = "323÷2=161, 1"; new = "407÷6=67, 5" },
    @{ old = "127÷3=42, 1"; new = "771÷5=154, 1" },
    @{ old = "317÷2=158, 1"; new = "697÷5=139, 2" },
    @{ old = "952÷3=317, 1"; new = "942÷2=471, 0" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
